$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet 1: Unit commitment
$ws1.Range("C2").Value = 13.9
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 13.9
$ws1.Range("F2").Value = 13.9
$ws1.Range("G2").Value = 22.92
$ws1.Range("H2").Value = 1509.48
$ws1.Range("I2").Value = 14866.68
$ws1.Range("J2").Value = 78.18000000000001
$ws1.Range("K2").Value = 769.73
$ws1.Range("L2").Value = 8.970000000000001
$ws1.Range("M2").Value = 84.22
$ws1.Range("N2").Value = 0.5
$ws1.Range("O2").Value = 9.6
$ws1.Range("P2").Value = 25.8
$ws1.Range("Q2").Value = 10.36
$ws1.Range("R2").Value = 3.3
$ws1.Range("T2").Value = 0
$ws1.Range("C3").Value = 13.98
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 13.98
$ws1.Range("F3").Value = 13.98
$ws1.Range("G3").Value = 22.58
$ws1.Range("H3").Value = 1511.63
$ws1.Range("I3").Value = 14866.68
$ws1.Range("J3").Value = 78.3
$ws1.Range("K3").Value = 769.73
$ws1.Range("L3").Value = 9.02
$ws1.Range("M3").Value = 84.22
$ws1.Range("N3").Value = 0.48
$ws1.Range("O3").Value = 9.6
$ws1.Range("P3").Value = 25.56
$ws1.Range("Q3").Value = 10.3
$ws1.Range("R3").Value = 3.06
$ws1.Range("T3").Value = 0
$ws1.Range("C4").Value = 13.65
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 13.65
$ws1.Range("F4").Value = 13.65
$ws1.Range("G4").Value = 22.06
$ws1.Range("H4").Value = 1502.69
$ws1.Range("I4").Value = 14866.68
$ws1.Range("J4").Value = 77.81999999999999
$ws1.Range("K4").Value = 769.73
$ws1.Range("L4").Value = 8.81
$ws1.Range("M4").Value = 84.22
$ws1.Range("N4").Value = 0.54
$ws1.Range("O4").Value = 9.6
$ws1.Range("P4").Value = 26.54
$ws1.Range("Q4").Value = 10.55
$ws1.Range("R4").Value = 4.04
$ws1.Range("T4").Value = 0
$ws1.Range("C5").Value = 13.36
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 13.36
$ws1.Range("F5").Value = 13.36
$ws1.Range("G5").Value = 21.58
$ws1.Range("H5").Value = 1494.7
$ws1.Range("I5").Value = 14866.68
$ws1.Range("J5").Value = 77.40000000000001
$ws1.Range("K5").Value = 769.73
$ws1.Range("L5").Value = 8.619999999999999
$ws1.Range("M5").Value = 84.22
$ws1.Range("N5").Value = 0.59
$ws1.Range("O5").Value = 9.6
$ws1.Range("P5").Value = 27.42
$ws1.Range("Q5").Value = 10.78
$ws1.Range("R5").Value = 4.92
$ws1.Range("T5").Value = 0
$ws1.Range("C6").Value = 13.14
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 13.14
$ws1.Range("F6").Value = 13.14
$ws1.Range("G6").Value = 21.23
$ws1.Range("H6").Value = 1488.81
$ws1.Range("I6").Value = 14866.68
$ws1.Range("J6").Value = 77.08
$ws1.Range("K6").Value = 769.73
$ws1.Range("L6").Value = 8.48
$ws1.Range("M6").Value = 84.22
$ws1.Range("N6").Value = 0.62
$ws1.Range("O6").Value = 9.6
$ws1.Range("P6").Value = 28.07
$ws1.Range("Q6").Value = 10.96
$ws1.Range("R6").Value = 5.57
$ws1.Range("T6").Value = 0
$ws1.Range("C7").Value = 12.92
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 12.92
$ws1.Range("F7").Value = 12.92
$ws1.Range("G7").Value = 20.86
$ws1.Range("H7").Value = 1482.71
$ws1.Range("I7").Value = 14866.68
$ws1.Range("J7").Value = 76.76000000000001
$ws1.Range("K7").Value = 769.73
$ws1.Range("L7").Value = 8.33
$ws1.Range("M7").Value = 84.22
$ws1.Range("N7").Value = 0.66
$ws1.Range("O7").Value = 9.6
$ws1.Range("P7").Value = 28.75
$ws1.Range("Q7").Value = 11.15
$ws1.Range("R7").Value = 6.25
$ws1.Range("T7").Value = 0
$ws1.Range("C8").Value = 12.69
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 12.69
$ws1.Range("F8").Value = 12.69
$ws1.Range("G8").Value = 20.5
$ws1.Range("H8").Value = 1476.65
$ws1.Range("I8").Value = 14866.68
$ws1.Range("J8").Value = 76.44
$ws1.Range("K8").Value = 769.73
$ws1.Range("L8").Value = 8.18
$ws1.Range("M8").Value = 84.22
$ws1.Range("N8").Value = 0.7
$ws1.Range("O8").Value = 9.6
$ws1.Range("P8").Value = 29.44
$ws1.Range("Q8").Value = 11.35
$ws1.Range("R8").Value = 6.94
$ws1.Range("T8").Value = 0
$ws1.Range("C9").Value = 12.46
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 12.46
$ws1.Range("F9").Value = 12.46
$ws1.Range("G9").Value = 20.13
$ws1.Range("H9").Value = 1470.64
$ws1.Range("I9").Value = 14866.68
$ws1.Range("J9").Value = 76.12
$ws1.Range("K9").Value = 769.73
$ws1.Range("L9").Value = 8.039999999999999
$ws1.Range("M9").Value = 84.22
$ws1.Range("N9").Value = 0.74
$ws1.Range("O9").Value = 9.6
$ws1.Range("P9").Value = 30.12
$ws1.Range("Q9").Value = 11.56
$ws1.Range("R9").Value = 7.62
$ws1.Range("T9").Value = 0
$ws1.Range("C10").Value = 12.23
$ws1.Range("D10").Value = 0
$ws1.Range("E10").Value = 12.23
$ws1.Range("F10").Value = 12.23
$ws1.Range("G10").Value = 19.76
$ws1.Range("H10").Value = 1464.67
$ws1.Range("I10").Value = 14866.68
$ws1.Range("J10").Value = 75.81
$ws1.Range("K10").Value = 769.73
$ws1.Range("L10").Value = 7.89
$ws1.Range("M10").Value = 84.22
$ws1.Range("N10").Value = 0.77
$ws1.Range("O10").Value = 9.6
$ws1.Range("P10").Value = 30.8
$ws1.Range("Q10").Value = 11.77
$ws1.Range("R10").Value = 8.300000000000001
$ws1.Range("T10").Value = 0
$ws1.Range("C11").Value = 12.23
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 12.23
$ws1.Range("F11").Value = 12.23
$ws1.Range("G11").Value = 19.33
$ws1.Range("H11").Value = 1464.7
$ws1.Range("I11").Value = 14866.68
$ws1.Range("J11").Value = 75.81
$ws1.Range("K11").Value = 769.73
$ws1.Range("L11").Value = 7.89
$ws1.Range("M11").Value = 84.22
$ws1.Range("O11").Value = 9.6
$ws1.Range("P11").Value = 30.8
$ws1.Range("Q11").Value = 11.77
$ws1.Range("R11").Value = 8.300000000000001
$ws1.Range("T11").Value = 0
$ws1.Range("C12").Value = 13.06
$ws1.Range("D12").Value = 0
$ws1.Range("E12").Value = 13.06
$ws1.Range("F12").Value = 13.06
$ws1.Range("G12").Value = 21.1
$ws1.Range("H12").Value = 1486.67
$ws1.Range("I12").Value = 14866.68
$ws1.Range("J12").Value = 76.97
$ws1.Range("K12").Value = 769.73
$ws1.Range("L12").Value = 8.42
$ws1.Range("M12").Value = 84.22
$ws1.Range("N12").Value = 0.64
$ws1.Range("O12").Value = 9.6
$ws1.Range("P12").Value = 28.33
$ws1.Range("Q12").Value = 11.05
$ws1.Range("R12").Value = 5.83
$ws1.Range("T12").Value = 0

# Sheet 2: Storages
$ws2.Range("C2").Value = 41.7
$ws2.Range("D2").Value = 22.92
$ws2.Range("E2").Value = 30.15
$ws2.Range("C3").Value = 41.94
$ws2.Range("D3").Value = 22.58
$ws2.Range("E3").Value = 35.32
$ws2.Range("C4").Value = 40.96
$ws2.Range("D4").Value = 22.06
$ws2.Range("E4").Value = 35.84
$ws2.Range("C5").Value = 40.08
$ws2.Range("D5").Value = 21.58
$ws2.Range("E5").Value = 36.32
$ws2.Range("C6").Value = 39.43
$ws2.Range("D6").Value = 21.23
$ws2.Range("E6").Value = 36.67
$ws2.Range("C7").Value = 38.75
$ws2.Range("D7").Value = 20.86
$ws2.Range("E7").Value = 37.04
$ws2.Range("C8").Value = 38.06
$ws2.Range("D8").Value = 20.5
$ws2.Range("E8").Value = 37.4
$ws2.Range("F8").Value = -0
$ws2.Range("G8").Value = 0.2
$ws2.Range("C9").Value = 37.38
$ws2.Range("D9").Value = 20.13
$ws2.Range("E9").Value = 30.13
$ws2.Range("G9").Value = 0.2
$ws2.Range("C10").Value = 36.7
$ws2.Range("D10").Value = 19.76
$ws2.Range("E10").Value = 26.88
$ws2.Range("F10").Value = -0
$ws2.Range("G10").Value = 0.2
$ws2.Range("C11").Value = 36.7
$ws2.Range("D11").Value = 19.33
$ws2.Range("E11").Value = 18.87
$ws2.Range("F11").Value = -0.8
$ws2.Range("C12").Value = 39.17
$ws2.Range("D12").Value = 21.1
$ws2.Range("E12").Value = 32.46
